$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "... dropping a key when an enemy is defeated, etc." ->
#    "... dropping a key when an enemy is defeated, since any enemy
#     could drop a key when defeated."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "dropping a key when an enemy is defeated, etc.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "dropping a key when an enemy is defeated, since any enemy could drop a key when defeated.",
    2) | Out-Null

# ------------------------------------------------------------------
# 2. Goon "Don't repeat yourself" paragraph rewrite
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "can be seen here as the Actor class has been inherited to create Enemy class. Goon class inherits from Enemy class. This ensures that code is reusable, not repeated and consistent in creating an object that has the same set of properties while having the freedom to extend the system.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can be seen here as the Goon class inherits from the Enemy class, while Enemy is an extension of the Actor class from the edu.monash.fit2099.engine package, resulting code to be not repeated, reusable and consistent in creating instances that has the same set of properties while having the freedom to extend the system.",
    2) | Out-Null

# ------------------------------------------------------------------
# 3. Ninja "Don't repeat yourself" paragraph rewrite
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "can be seen here too as the Actor class has been inherited to create Enemy class. Ninja class inherits from Enemy class. This ensures that code is reusable, not repeated and consistent in creating an object that has the same set of properties while having the freedom to extend the system.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can be seen here too as the Ninja class inherits from the Enemy class, while Enemy is an extension of the Actor class from the edu.monash.fit2099.engine package, resulting code to be not repeated, reusable and consistent in creating instances that has the same set of properties while having the freedom to extend the system.",
    2) | Out-Null

# ------------------------------------------------------------------
# 4. Move the "_GoBack" bookmark so that it sits right after the
#    rewritten Goon paragraph (where the author's cursor last was),
#    instead of its original spot near the "Miniboss" heading.
#
#    Directly creating a zero-length bookmark at that exact position
#    can misplace it, so we anchor it with a temporary marker, insert
#    the bookmark next to the marker, and then remove the marker text
#    - the bookmark stays put because it is anchored to the
#    surrounding text, not the raw character offset.
# ------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Goon class inherits from the Enemy class*") {
        $insertPos = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        $insertPos.InsertAfter("ZZMARKERZZ")
        break
    }
}

$r = $d.Content
$r.Find.Execute("ZZMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$markerPos = $r.Start
$d.Bookmarks.Add("_GoBack", $d.Range($markerPos, $markerPos)) | Out-Null

$d.Content.Find.Execute("ZZMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
